$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header row cell text. Order matters: it determines the order new
# strings are appended to the shared strings table (B1, then C1, then A1
# matches the target table ordering: UNIDAD DE MEDIDA, SKU, CONCEPTO MATERIAL).
$ws.Range("B1").Value = "UNIDAD DE MEDIDA"
$ws.Range("C1").Value = "SKU"
$ws.Range("A1").Value = " CONCEPTO MATERIAL"

# Widen column B to fit the new, longer header text
$ws.Columns.Item(2).ColumnWidth = 17.87

# Update the active selection to A10
$ws.Range("A10").Select()
